$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read existing data block (rows 2-49, columns A-I) into memory.
# (Row 2 is the first data row "2014-01"; row 49 is the last "2017-12".)
$src = $ws.Range("A2:I49").Value()

$nRows = $src.GetLength(0)
$nCols = $src.GetLength(1)

# Build destination row order: within each block of 12 consecutive rows
# (one calendar year, originally ordered Jan..Dec), move the last 3 rows
# (Oct, Nov, Dec) to the front of that block, keeping Jan..Sep after them.
$order = New-Object 'object[]' $nRows
$blockSize = 12
$nBlocks = [int]($nRows / $blockSize)
for ($b = 0; $b -lt $nBlocks; $b++) {
    $base = $b * $blockSize
    $pos = $base
    for ($k = 9; $k -lt 12; $k++) {
        $order[$pos] = $base + $k
        $pos = $pos + 1
    }
    for ($k = 0; $k -lt 9; $k++) {
        $order[$pos] = $base + $k
        $pos = $pos + 1
    }
}

# Build the new 0-based value array applying the computed order.
$dest = New-Object 'object[,]' $nRows, $nCols
for ($r = 0; $r -lt $nRows; $r++) {
    $srcRow = [int]$order[$r]
    for ($c = 0; $c -lt $nCols; $c++) {
        # $src is 1-based (Range.Value semantics); $dest is 0-based.
        $dest[$r, $c] = $src[$srcRow + 1, $c + 1]
    }
}

$ws.Range("A2:I49").Value = $dest
